$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Localize header row labels (English -> Japanese) for the regional columns
# that received translations in this revision. "Mountain" and "South" were
# left as-is in the source edit, so they are untouched here.
$ws.Range("A1").Value = "年 - 四半期"
$ws.Range("B1").Value = "中西部"
$ws.Range("D1").Value = "北東部"
$ws.Range("F1").Value = "東南アジア"
$ws.Range("G1").Value = "西部"
